{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the inside of: async (context) => { ... }\n//\n// Replaces the date line and the 25 \"three-digit \u00f7 one-digit\" answer\n// cells with their new values, per the commit diff. All old values are\n// unique substrings in the document, so a body.search() + set-text\n// approach is safe (no accidental double-replacement).\n\nconst replacements = [\n  [\"2024-09-14 Saturday\", \"2024-09-15 Sunday\"],\n  [\"237\u00f72=118, 1\", \"502\u00f75=100, 2\"],\n  [\"249\u00f72=124, 1\", \"669\u00f75=133, 4\"],\n  [\"679\u00f73=226, 1\", \"420\u00f74=105, 0\"],\n  [\"730\u00f72=365, 0\", \"875\u00f74=218, 3\"],\n  [\"893\u00f74=223, 1\", \"590\u00f76=98, 2\"],\n  [\"252\u00f75=50, 2\", \"585\u00f77=83, 4\"],\n  [\"235\u00f78=29, 3\", \"211\u00f76=35, 1\"],\n  [\"995\u00f75=199, 0\", \"554\u00f74=138, 2\"],\n  [\"965\u00f79=107, 2\", \"381\u00f74=95, 1\"],\n  [\"268\u00f77=38, 2\", \"979\u00f79=108, 7\"],\n  [\"335\u00f74=83, 3\", \"106\u00f75=21, 1\"],\n  [\"264\u00f76=44, 0\", \"631\u00f73=210, 1\"],\n  [\"719\u00f78=89, 7\", \"982\u00f77=140, 2\"],\n  [\"428\u00f79=47, 5\", \"937\u00f74=234, 1\"],\n  [\"601\u00f72=300, 1\", \"645\u00f72=322, 1\"],\n  [\"583\u00f73=194, 1\", \"888\u00f78=111, 0\"],\n  [\"988\u00f72=494, 0\", \"274\u00f74=68, 2\"],\n  [\"243\u00f79=27, 0\", \"370\u00f79=41, 1\"],\n  [\"492\u00f74=123, 0\", \"923\u00f78=115, 3\"],\n  [\"290\u00f78=36, 2\", \"773\u00f78=96, 5\"],\n  [\"594\u00f72=297, 0\", \"991\u00f75=198, 1\"],\n  [\"446\u00f75=89, 1\", \"533\u00f77=76, 1\"],\n  [\"869\u00f78=108, 5\", \"536\u00f77=76, 4\"],\n  [\"159\u00f79=17, 6\", \"525\u00f75=105, 0\"],\n  [\"948\u00f75=189, 3\", \"949\u00f76=158, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script.\n# $word / $doc / $app resolve to the running Word session; the active\n# document is $word.ActiveDocument.\n#\n# Replaces the date line and the 25 \"three-digit \u00f7 one-digit\" answer\n# cells with their new values, per the commit diff, using Find/Replace\n# (wdReplaceAll) on each unique old string.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-09-14 Saturday\", \"2024-09-15 Sunday\"),\n    @(\"237\u00f72=118, 1\", \"502\u00f75=100, 2\"),\n    @(\"249\u00f72=124, 1\", \"669\u00f75=133, 4\"),\n    @(\"679\u00f73=226, 1\", \"420\u00f74=105, 0\"),\n    @(\"730\u00f72=365, 0\", \"875\u00f74=218, 3\"),\n    @(\"893\u00f74=223, 1\", \"590\u00f76=98, 2\"),\n    @(\"252\u00f75=50, 2\", \"585\u00f77=83, 4\"),\n    @(\"235\u00f78=29, 3\", \"211\u00f76=35, 1\"),\n    @(\"995\u00f75=199, 0\", \"554\u00f74=138, 2\"),\n    @(\"965\u00f79=107, 2\", \"381\u00f74=95, 1\"),\n    @(\"268\u00f77=38, 2\", \"979\u00f79=108, 7\"),\n    @(\"335\u00f74=83, 3\", \"106\u00f75=21, 1\"),\n    @(\"264\u00f76=44, 0\", \"631\u00f73=210, 1\"),\n    @(\"719\u00f78=89, 7\", \"982\u00f77=140, 2\"),\n    @(\"428\u00f79=47, 5\", \"937\u00f74=234, 1\"),\n    @(\"601\u00f72=300, 1\", \"645\u00f72=322, 1\"),\n    @(\"583\u00f73=194, 1\", \"888\u00f78=111, 0\"),\n    @(\"988\u00f72=494, 0\", \"274\u00f74=68, 2\"),\n    @(\"243\u00f79=27, 0\", \"370\u00f79=41, 1\"),\n    @(\"492\u00f74=123, 0\", \"923\u00f78=115, 3\"),\n    @(\"290\u00f78=36, 2\", \"773\u00f78=96, 5\"),\n    @(\"594\u00f72=297, 0\", \"991\u00f75=198, 1\"),\n    @(\"446\u00f75=89, 1\", \"533\u00f77=76, 1\"),\n    @(\"869\u00f78=108, 5\", \"536\u00f77=76, 4\"),\n    @(\"159\u00f79=17, 6\", \"525\u00f75=105, 0\"),\n    @(\"948\u00f75=189, 3\", \"949\u00f76=158, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2) | Out-Null\n}\n"}
